# edit.ps1
# Applies the JOSS paper v0.6 text revisions to the ggoncoplot paper.
#
# Notes on technique:
#  - A plain Find/Replace on $d.Content, when it touches a paragraph,
#    causes this runtime to coalesce *all* identically-formatted runs in
#    that paragraph on save (even runs untouched by the actual edit).
#    The target document keeps several pre-existing same-format run
#    splits (e.g. around inline citation markers) that must NOT be
#    touched by edits elsewhere in the same paragraph. Doing the
#    replace with TrackRevisions on, then immediately calling
#    Revisions.AcceptAll(), avoids that unwanted coalescing: only the
#    literal changed span becomes (after accepting) a single new run,
#    and sibling runs elsewhere in the paragraph are left exactly as
#    they were.
#  - To match this behaviour exactly, each Find/Replace below always
#    targets the *entire* text of the run(s) being changed (not just a
#    sub-string), so that after accepting, the edited span collapses
#    back into one run rather than being split at the edit boundary.
#  - One bullet ("Interactive plots") goes from one run to three runs
#    in the target, because a brand-new sentence is appended and then
#    two further (genuinely new) runs follow it (" " and "."),
#    mirroring a reference/citation pattern already used elsewhere in
#    this document. New sibling runs like that are created reliably by
#    calling InsertAfter on a freshly resolved zero-length Range
#    (re-fetching Range.End between inserts), rather than reusing /
#    duplicating a Range handle.

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.TrackRevisions = $true
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        $d.TrackRevisions = $false
        throw "Find failed for: $oldText"
    }
    $d.TrackRevisions = $false
    [void]$d.Revisions.AcceptAll()
}

# 1 & 2) Summary paragraph, first two runs.
Replace-Text `
    "The ggoncoplot R package generates interactive oncoplots (also called oncoprints) to visualize mutational patterns across patient cancer cohorts (" `
    "The ggoncoplot R package generates interactive oncoplots to visualize mutational patterns across patient cancer cohorts ("

Replace-Text `
    "). Oncoplots reveal patterns of gene co-mutation and include marginal plots that indicate co-occurrence of gene mutations and tumour features. It is useful to relate gene mutation patterns seen in an oncoplot to patterns seen in other plot types, including gene expression t-SNE plots or methylation UMAPs. There are, however, no existing oncoplot-generating R packages that support dynamic data linkage between different plots. To addresses this gap and enable rapid exploration of a variety of data types we constructed the ggoncoplot package for the production of oncoplots that are easily integrated with custom visualisations and that support synchronised data-selections across plots (" `
    "). Oncoplots, also called oncoprints, reveal patterns of gene co-mutation and include marginal plots that indicate co-occurrence of gene mutations with tumour and clinical features. It is useful to relate gene mutation patterns seen in an oncoplot to patterns in other plot types, including gene expression t-SNE plots or methylation UMAPs. The simplest and most intuitive approach to examining such relations is to link plots dynamically such that samples selected in an oncoplot can be highlighted in other plots. There are, however, no existing oncoplot-generating R packages that support dynamic data linkage between different plots. To address this gap and enable rapid exploration of a variety of data types we constructed the ggoncoplot package for the production of oncoplots that are easily integrated with custom visualisations and that support synchronised data-selections across plots ("

# 3) Statement of Need, first run.
Replace-Text `
    "Oncoplots are highly effectively for visualising mutation data in cancer cohorts but are challenging to generate with the major R plotting systems (base, lattice, or ggplot2) due to their algorithmic and graphical complexity. Simplifying the process would make oncoplots more accessible to researchers. Packages like ComplexHeatmap" `
    "Oncoplots are highly effective for visualising mutation data in cancer cohorts but are challenging to generate with the major R plotting systems (base, lattice, or ggplot2) due to their algorithmic and graphical complexity. Simplifying the process of generating oncoplots would make them more accessible to researchers. Existing packages including ComplexHeatmap"

# 4) "all make static oncoplots..." run.
Replace-Text `
    "all make static oncoplots easier to create, but there is still a significant unmet need for an easy method of creating oncoplots with the following features:" `
    "all make static oncoplots easier to create, but there is still a significant unmet need for user-friendly method of creating oncoplots with the following features:"

# 5) "Interactive plots" bullet: extend the existing run, then append two
#    brand-new runs (" " and ".") after it, matching the diff's 3-run
#    structure.
Replace-Text `
    ": Customizable tooltips, cross-selection of samples across different plots, and auto-copying of sample identifiers on click." `
    ": Customizable tooltips, cross-selection of samples across different plots, and auto-copying of sample identifiers on click. This enables exploration of trends in multiomic datasets as shown in"

$hit = $d.Content
$found = $hit.Find.Execute("as shown in", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Find failed for: as shown in"
}
$insertPos = $hit.End
$rSpace = $d.Range($insertPos, $insertPos)
$rSpace.InsertAfter(" ")
$rPeriod = $d.Range($rSpace.End, $rSpace.End)
$rPeriod.InsertAfter(".")

# 6) "Support for tidy datasets" bullet.
Replace-Text `
    ": Compatibility with tidy, tabular mutation-level formats (MAF files or relational databases), typical of cancer cohort datasets." `
    ": Compatibility with tidy, tabular mutation-level formats (MAF files or relational databases), typical of cancer cohort datasets. This greatly improves the range of datasets that can be quickly and easily visualised in an oncoplot."

# 7) "Auto colouring" -> "Auto-colouring" (bold heading run).
Replace-Text "Auto colouring" "Auto-colouring"

# 8) Auto-colouring bullet body.
Replace-Text `
    ": Automatic selection of colour palettes for datasets where consequence annotations are aligned with standard variant effect dictionaries (PAVE, SO, or MAF)." `
    ": Automatic selection of colour palettes for datasets where the consequence annotations are aligned with standard variant effect dictionaries (PAVE, SO, or MAF)."

# 9) "Versatility" bullet body.
Replace-Text `
    ": The ability to visualize entities other than gene mutations, including noncoding features (e.g., enhancers) and non-genomic entities (e.g., microbial presence in microbiome datasets)." `
    ": The ability to visualize entities other than gene mutations, such as noncoding features (e.g., promoter or enhancer mutations) and non-genomic entities (e.g., microbial presence in microbiome datasets)."

# 10) "We developed ggoncoplot..." run.
Replace-Text `
    "We developed ggoncoplot as the first R package that addresses all these challenges simultaneously (" `
    "We developed ggoncoplot as the first R package to address all these challenges together ("
